# Applies the operator-precedence-table edit described by the diff:
#  - Fills the two previously-empty cells in the "+" row (table row 2)
#    with bold 20pt "<" characters.
#  - Fills the two previously-empty cells in the "-" row (table row 3)
#    with bold 20pt "<" characters, and relocates the document's
#    "_GoBack" bookmark onto the second of those cells.
#  - Fills the two previously-empty cells in the "*" row (table row 4)
#    with bold 20pt ">" characters.
#  - Fills the two previously-empty cells in the "/" row (table row 5)
#    with bold 20pt ">" characters.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Fill-Cell($table, $rowIdx, $colIdx, $char) {
    $cell = $table.Cell($rowIdx, $colIdx)
    $rng = $cell.Range
    $rng.InsertBefore($char)

    # Re-fetch the cell - the handle becomes stale once the table text changes.
    $cell2 = $table.Cell($rowIdx, $colIdx)
    $rng2 = $cell2.Range
    $rng2.Font.Bold = 1
    $rng2.Font.Size = 20
    $rng2.Font.SizeBi = 20
}

# Row 2 ("+"): columns 12 and 13 get "<"
Fill-Cell $t 2 12 "<"
Fill-Cell $t 2 13 "<"

# Row 3 ("-"): columns 12 and 13 get "<"
Fill-Cell $t 3 12 "<"
Fill-Cell $t 3 13 "<"

# Row 4 ("*"): columns 12 and 13 get ">"
Fill-Cell $t 4 12 ">"
Fill-Cell $t 4 13 ">"

# Row 5 ("/"): columns 12 and 13 get ">"
Fill-Cell $t 5 12 ">"
Fill-Cell $t 5 13 ">"

# Relocate the "_GoBack" bookmark from its old spot (row 13, column 5) to
# the newly-filled cell at row 3, column 13, matching the diff.
$targetCell = $t.Cell(3, 13)
$targetRange = $targetCell.Range
# Collapse to just before the cell's paragraph mark (i.e. right after the
# "<" run that was just inserted), mirroring where the bookmark sits in
# the target XML.
$bmRange = $d.Range($targetRange.Start, $targetRange.End - 1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)
